$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: Volume/Number and report week date range ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Crime Complaints table (Week to Date / 28 Day / YTD / 2 Year) updates ---
# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 7.692307692307
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 117
$ws.Range("K16").Value = -15.384615384615
$ws.Range("L16").Value = -18.181818181818
$ws.Range("M16").Value = -37.735849056603
# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 31.818181818181
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 181
$ws.Range("K17").Value = 10.497237569060
$ws.Range("L17").Value = 9.890109890109
$ws.Range("M17").Value = 55.038759689922
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -7.692307692307
$ws.Range("I18").Value = 74
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = -5.128205128205
$ws.Range("L18").Value = -13.953488372093
$ws.Range("M18").Value = 4.225352112676
# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 125
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -15.625
$ws.Range("I19").Value = 262
$ws.Range("J19").Value = 262
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2.34375
$ws.Range("M19").Value = 42.391304347826
# Row 20
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 104
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 35.064935064935
$ws.Range("L20").Value = -0.952380952380
$ws.Range("M20").Value = 121.276595744681
# Row 21
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 15.789473684210
$ws.Range("F21").Value = 94
$ws.Range("H21").Value = 6.818181818181
$ws.Range("I21").Value = 753
$ws.Range("J21").Value = 726
$ws.Range("K21").Value = 3.719008264462
$ws.Range("L21").Value = -1.181102362204
$ws.Range("M21").Value = 24.875621890547
# Row 22
$ws.Range("D22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -55.555555555555
# Row 23
$ws.Range("C23").Value = 1
$ws.Range("F23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 11
$ws.Range("K23").Value = -47.619047619047
$ws.Range("L23").Value = -8.333333333333
$ws.Range("M23").Value = 22.222222222222
# Row 24
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -13.333333333333
$ws.Range("F24").Value = 74
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 13.846153846153
$ws.Range("I24").Value = 569
$ws.Range("J24").Value = 640
$ws.Range("K24").Value = -11.09375
$ws.Range("L24").Value = -16.934306569343
$ws.Range("M24").Value = 41.191066997518
# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 15
$ws.Range("H25").Value = 6.666666666666
$ws.Range("I25").Value = 119
$ws.Range("J25").Value = 202
$ws.Range("K25").Value = -41.089108910891
$ws.Range("L25").Value = -28.313253012048
# Row 26
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 27
$ws.Range("H26").Value = 28.571428571428
$ws.Range("I26").Value = 266
$ws.Range("J26").Value = 236
$ws.Range("K26").Value = 12.711864406779
$ws.Range("L26").Value = -6.007067137809
$ws.Range("M26").Value = -22.674418604651
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("F27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 15
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 15.384615384615
# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 300
$ws.Range("L28").Value = 40

$excel.CutCopyMode = $false
